$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Gori"

# Mark confidential/unavailable data points with "..."
$ws.Range("H6").Value = "..."
$ws.Range("O6").Value = "..."
$ws.Range("O7").Value = "..."

# Add footnote row explaining the "..." marker
$ws.Range("A8").Value = "Note: „ ... „ - Data is confidential or unavailable."
$noteChars = $ws.Range("A8").Characters(1, 5)
$noteChars.Font.Bold = $true
$noteChars.Font.Underline = $true
$noteChars.Font.Size = 9
$noteChars.Font.Name = "Arial"
$restChars = $ws.Range("A8").Characters(6, 47)
$restChars.Font.Size = 9
$restChars.Font.Name = "Arial"
$ws.Range("A8").Font.Size = 9

Write-Host "done"
